# Apply cryptos list update (Tue Aug 20 14:50:53 UTC 2024) via GitHub Actions
# Updates price (column D) and volume-change (column E) text values for
# existing rows, and swaps the WhiteBITCoin/Hedera rows (46/47) with new data.
# Column D values are forced to Text format first so numeric-looking prices
# (e.g. "1.00", "59.552.23") are preserved exactly as text, matching the
# source data feed's formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.552.23"
$ws.Range("E2").Value = "  +2.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.592.29"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.02"
$ws.Range("E5").Value = "  +4.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.20"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.596"
$ws.Range("E8").Value = "  +2.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.600.34"
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("E10").Value = "  -1.73%  "
$ws.Range("E11").Value = "  +3.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.151"
$ws.Range("E12").Value = "  +9.50%  "
$ws.Range("E13").Value = "  +2.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.049.86"
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.553.17"
$ws.Range("E15").Value = "  +2.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.88"
$ws.Range("E16").Value = "  +6.56%  "
$ws.Range("E17").Value = "  +3.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.599.73"
$ws.Range("E18").Value = "  +1.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.51"
$ws.Range("E19").Value = "  +1.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "336.63"
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.20"
$ws.Range("E21").Value = "  +2.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.24"
$ws.Range("E22").Value = "  +2.23%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.93"
$ws.Range("E24").Value = "  -2.07%  "
$ws.Range("E25").Value = "  +5.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  +2.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.26"
$ws.Range("E28").Value = "  +3.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0780"
$ws.Range("E29").Value = "  +6.80%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.69"
$ws.Range("E31").Value = "  +3.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "160.35"
$ws.Range("E32").Value = "  +4.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.03"
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.95"
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.03"
$ws.Range("E35").Value = "  +3.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.889"
$ws.Range("E36").Value = "  +9.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.877"
$ws.Range("E37").Value = "  +3.29%  "
$ws.Range("E38").Value = "  +4.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.96"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("E40").Value = "  +5.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "292.23"
$ws.Range("E41").Value = "  +4.80%  "
$ws.Range("E42").Value = "  +1.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.997"
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0974"
$ws.Range("E44").Value = "  +3.61%  "
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.66"
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0536"
$ws.Range("E47").Value = "  +1.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.96"
$ws.Range("E48").Value = "  +3.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.92"
$ws.Range("E49").Value = "  +14.91%  "
$ws.Range("E50").Value = "  +2.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.931.23"
$ws.Range("E51").Value = "  +1.60%  "
